# Fix wrong concept for onTouchEvent
# Swap the order of the "Interested View Example" and "Ignorant View Example"
# slides (currently slides 7 and 8).
$p = $ppt.ActivePresentation

$s = $p.Slides.Item(7)
$s.MoveTo(8)
